# Update "F" column figures (e.g. view/attendance counts) across sheets
# to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 327
$ws1.Range("F3").Value = 3335
$ws1.Range("F6").Value = 7794
$ws1.Range("F14").Value = 1783
$ws1.Range("F15").Value = 387
$ws1.Range("F16").Value = 119
$ws1.Range("F22").Value = 8067
$ws1.Range("F33").Value = 536
$ws1.Range("F34").Value = 536
$ws1.Range("F42").Value = 3328
$ws1.Range("F44").Value = 728
$ws1.Range("F48").Value = 320

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 306
$ws2.Range("F19").Value = 46

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 2179
$ws3.Range("F10").Value = 1162

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 327
$ws4.Range("F5").Value = 7794
$ws4.Range("F9").Value = 1162
$ws4.Range("F15").Value = 306
$ws4.Range("F16").Value = 387
$ws4.Range("F17").Value = 119
$ws4.Range("F22").Value = 8067
$ws4.Range("F30").Value = 536
$ws4.Range("F31").Value = 536
$ws4.Range("F39").Value = 3328
